# Remove type_lieu from contrat modal
# The "Périodicité" value that used to live (mistakenly) in column A now
# gets its own column (F). Row 2 is updated with the corrected contract
# data, and a second contract (row 3) is appended below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : update existing contract row ---------------------------------
$ws.Range("A2").Value = "444/AA4444"
$ws.Range("B2").Value = "Point de vente"

# C2 looks numeric ("101") but must stay text, like the rest of the column.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "101"

$ws.Range("D2").Value = "mmmmml"
$ws.Range("E2").Value = "ds"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 12200
$ws.Range("H2").Value = 12200
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 1830
$ws.Range("K2").Value = 1830
$ws.Range("L2").Value = 24400
$ws.Range("M2").Value = 10370

# --- Row 3 : new contract row ----------------------------------------------
$ws.Range("A3").Value = "006/tESTDRR"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "aaaaa"
$ws.Range("D3").Value = "aaaaaaa"
$ws.Range("E3").Value = "ds"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 10000
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 8500
